$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 0
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 14).ClearContents()

$ws.Cells.Item(51, 8).Value = 7000
$ws.Cells.Item(51, 10).Value = 7000
$ws.Cells.Item(51, 12).Value = 7000
$ws.Cells.Item(51, 14).Value = -7968

$ws.Cells.Item(53, 8).Value = 7397.6875
$ws.Cells.Item(53, 9).Value = 9004.846
$ws.Cells.Item(53, 10).Value = 433.33334
$ws.Cells.Item(53, 11).Value = 9004.846
$ws.Cells.Item(53, 12).Value = 433.33334
$ws.Cells.Item(53, 13).Value = -8367.846
$ws.Cells.Item(53, 14).Value = -1707.33334

$ws.Cells.Item(62, 8).Value = 1285
$ws.Cells.Item(62, 9).Value = 1332.5
$ws.Cells.Item(62, 11).Value = 1332.5
$ws.Cells.Item(62, 13).Value = -708.5

$ws.Cells.Item(65, 8).Value = 1285
$ws.Cells.Item(65, 9).Value = 1332.5
$ws.Cells.Item(65, 11).Value = 6662.5
$ws.Cells.Item(65, 13).Value = -3542.5

$ws.Cells.Item(125, 8).Value = 1508.3334
$ws.Cells.Item(125, 9).Value = 1488.1666
$ws.Cells.Item(125, 10).Value = 1528.5
$ws.Cells.Item(125, 11).Value = 13393.4994
$ws.Cells.Item(125, 12).Value = 13756.5
$ws.Cells.Item(125, 13).Value = -10933.4994
$ws.Cells.Item(125, 14).Value = -18676.5

$ws.Cells.Item(129, 8).Value = 883.1163
$ws.Cells.Item(129, 10).Value = 892
$ws.Cells.Item(129, 12).Value = 2676
$ws.Cells.Item(129, 14).Value = -12676

$ws.Cells.Item(132, 8).Value = 575.4691
$ws.Cells.Item(132, 9).Value = 514.7895
$ws.Cells.Item(132, 11).Value = 1544.3685
$ws.Cells.Item(132, 13).Value = 985.6315

$ws.Cells.Item(137, 8).Value = 2339.3333
$ws.Cells.Item(137, 9).Value = 1733
$ws.Cells.Item(137, 10).Value = 2460.6
$ws.Cells.Item(137, 11).Value = 5199
$ws.Cells.Item(137, 12).Value = 7381.799999999999
$ws.Cells.Item(137, 13).Value = -2649
$ws.Cells.Item(137, 14).Value = -12481.8

$ws.Cells.Item(138, 8).Value = 2418
$ws.Cells.Item(138, 9).Value = 2515.926
$ws.Cells.Item(138, 10).Value = 2316.3076
$ws.Cells.Item(138, 11).Value = 7547.778
$ws.Cells.Item(138, 12).Value = 6948.9228
$ws.Cells.Item(138, 13).Value = -2407.778
$ws.Cells.Item(138, 14).Value = -17228.9228

$ws.Cells.Item(140, 8).Value = 53748.047
$ws.Cells.Item(140, 10).Value = 53748.047
$ws.Cells.Item(140, 12).Value = 53748.047
$ws.Cells.Item(140, 14).Value = -64108.047

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2814.348
$ws.Cells.Item(32, 9).Value = 2156.6309
$ws.Cells.Item(32, 11).Value = 2156.6309
$ws.Cells.Item(32, 13).Value = -1869.6309

$ws.Cells.Item(46, 8).Value = 11337.429
$ws.Cells.Item(46, 9).Value = 10137
$ws.Cells.Item(46, 10).Value = 12237.75
$ws.Cells.Item(46, 11).Value = 10137
$ws.Cells.Item(46, 12).Value = 12237.75
$ws.Cells.Item(46, 13).Value = -9818
$ws.Cells.Item(46, 14).Value = -12875.75

$ws.Cells.Item(61, 8).Value = 3137.5557
$ws.Cells.Item(61, 9).Value = 1654.75
$ws.Cells.Item(61, 11).Value = 1654.75
$ws.Cells.Item(61, 13).Value = -1442.75

$ws.Cells.Item(74, 8).Value = 2292.6667
$ws.Cells.Item(74, 9).Value = 1495.6666
$ws.Cells.Item(74, 11).Value = 1495.6666
$ws.Cells.Item(74, 13).Value = -621.6666

$ws.Cells.Item(77, 8).Value = 2292.6667
$ws.Cells.Item(77, 9).Value = 1495.6666
$ws.Cells.Item(77, 11).Value = 7478.333000000001
$ws.Cells.Item(77, 13).Value = -3110.333000000001

$ws.Cells.Item(102, 8).Value = 2299.4
$ws.Cells.Item(102, 9).Value = 2299.4
$ws.Cells.Item(102, 11).Value = 2299.4
$ws.Cells.Item(102, 13).Value = -677.4000000000001

$ws.Cells.Item(110, 8).Value = 1872
$ws.Cells.Item(110, 9).Value = 1601.4546
$ws.Cells.Item(110, 10).Value = 3062.4
$ws.Cells.Item(110, 11).Value = 1601.4546
$ws.Cells.Item(110, 12).Value = 3062.4
$ws.Cells.Item(110, 13).Value = 443.5454
$ws.Cells.Item(110, 14).Value = -7152.4

$ws.Cells.Item(122, 8).Value = 1106.9546
$ws.Cells.Item(122, 9).Value = 1243.9333
$ws.Cells.Item(122, 11).Value = 3731.7999
$ws.Cells.Item(122, 13).Value = -1281.7999

$ws.Cells.Item(132, 8).Value = 1901.6364
$ws.Cells.Item(132, 9).Value = 1319.9412
$ws.Cells.Item(132, 10).Value = 3879.4
$ws.Cells.Item(132, 11).Value = 3959.8236
$ws.Cells.Item(132, 12).Value = 11638.2
$ws.Cells.Item(132, 13).Value = -1429.8236
$ws.Cells.Item(132, 14).Value = -16698.2

$ws.Cells.Item(136, 8).Value = 3137.5557
$ws.Cells.Item(136, 9).Value = 1654.75
$ws.Cells.Item(136, 11).Value = 4964.25
$ws.Cells.Item(136, 13).Value = -2414.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(34, 8).Value = 10000
$ws.Cells.Item(34, 10).Value = 10000
$ws.Cells.Item(34, 12).Value = 10000
$ws.Cells.Item(34, 14).Value = -10228

$ws.Cells.Item(81, 8).Value = 22509.666
$ws.Cells.Item(81, 10).Value = 22509.666
$ws.Cells.Item(81, 12).Value = 22509.666
$ws.Cells.Item(81, 14).Value = -24631.666

$ws.Cells.Item(84, 8).Value = 22509.666
$ws.Cells.Item(84, 10).Value = 22509.666
$ws.Cells.Item(84, 12).Value = 67528.99800000001
$ws.Cells.Item(84, 14).Value = -78136.99800000001

$ws.Cells.Item(134, 8).Value = 5803
$ws.Cells.Item(134, 9).Value = 7118.52
$ws.Cells.Item(134, 10).Value = 2813.182
$ws.Cells.Item(134, 11).Value = 21355.56
$ws.Cells.Item(134, 12).Value = 8439.545999999998
$ws.Cells.Item(134, 13).Value = -18820.56
$ws.Cells.Item(134, 14).Value = -13509.546

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1127.6154
$ws.Cells.Item(22, 9).Value = 300.5
$ws.Cells.Item(22, 10).Value = 1278
$ws.Cells.Item(22, 11).Value = 300.5
$ws.Cells.Item(22, 12).Value = 1278
$ws.Cells.Item(22, 13).Value = 49.5
$ws.Cells.Item(22, 14).Value = -1978

$ws.Cells.Item(31, 8).Value = 3151.3572
$ws.Cells.Item(31, 9).Value = 3741.4
$ws.Cells.Item(31, 10).Value = 2823.5557
$ws.Cells.Item(31, 11).Value = 3741.4
$ws.Cells.Item(31, 12).Value = 2823.5557
$ws.Cells.Item(31, 13).Value = -3446.4
$ws.Cells.Item(31, 14).Value = -3413.5557

$ws.Cells.Item(34, 8).Value = 3151.3572
$ws.Cells.Item(34, 9).Value = 3741.4
$ws.Cells.Item(34, 10).Value = 2823.5557
$ws.Cells.Item(34, 11).Value = 3741.4
$ws.Cells.Item(34, 12).Value = 2823.5557
$ws.Cells.Item(34, 13).Value = -3539.4
$ws.Cells.Item(34, 14).Value = -3227.5557

$ws.Cells.Item(58, 8).Value = 1209628.8
$ws.Cells.Item(58, 9).Value = 1673604.2
$ws.Cells.Item(58, 11).Value = 1673604.2
$ws.Cells.Item(58, 13).Value = -1673401.2

$ws.Cells.Item(122, 8).Value = 3589.6365
$ws.Cells.Item(122, 9).Value = 995.6667
$ws.Cells.Item(122, 10).Value = 6702.4
$ws.Cells.Item(122, 11).Value = 2987.0001
$ws.Cells.Item(122, 12).Value = 20107.2
$ws.Cells.Item(122, 13).Value = -537.0001000000002
$ws.Cells.Item(122, 14).Value = -25007.2

$ws.Cells.Item(132, 8).Value = 2126.5293
$ws.Cells.Item(132, 9).Value = 1147.9524
$ws.Cells.Item(132, 11).Value = 3443.857199999999
$ws.Cells.Item(132, 13).Value = -913.8571999999995

$ws.Cells.Item(134, 8).Value = 2169.1155
$ws.Cells.Item(134, 10).Value = 4333
$ws.Cells.Item(134, 12).Value = 12999
$ws.Cells.Item(134, 14).Value = -18069

$ws.Cells.Item(136, 8).Value = 1209628.8
$ws.Cells.Item(136, 9).Value = 1673604.2
$ws.Cells.Item(136, 11).Value = 5020812.6
$ws.Cells.Item(136, 13).Value = -5018262.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 6421303
$ws.Cells.Item(131, 9).Value = 100000420
$ws.Cells.Item(131, 10).Value = 11774.863
$ws.Cells.Item(131, 11).Value = 300001260
$ws.Cells.Item(131, 12).Value = 35324.589
$ws.Cells.Item(131, 13).Value = -299996220
$ws.Cells.Item(131, 14).Value = -45404.589

$ws.Cells.Item(132, 8).Value = 1082.6
$ws.Cells.Item(132, 9).Value = 690
$ws.Cells.Item(132, 10).Value = 1126.2222
$ws.Cells.Item(132, 11).Value = 6210
$ws.Cells.Item(132, 12).Value = 10135.9998
$ws.Cells.Item(132, 13).Value = -3680
$ws.Cells.Item(132, 14).Value = -15195.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2027758.5
$ws.Cells.Item(132, 9).Value = 3848523.5
$ws.Cells.Item(132, 10).Value = 4686.3335
$ws.Cells.Item(132, 11).Value = 11545570.5
$ws.Cells.Item(132, 12).Value = 14059.0005
$ws.Cells.Item(132, 13).Value = -11543040.5
$ws.Cells.Item(132, 14).Value = -19119.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 4753.75
$ws.Cells.Item(16, 10).Value = 4102
$ws.Cells.Item(16, 12).Value = 4102
$ws.Cells.Item(16, 14).Value = -4442

$ws.Cells.Item(22, 8).Value = 2416.6
$ws.Cells.Item(22, 9).Value = 2691.6667
$ws.Cells.Item(22, 10).Value = 2233.2222
$ws.Cells.Item(22, 11).Value = 2691.6667
$ws.Cells.Item(22, 12).Value = 2233.2222
$ws.Cells.Item(22, 13).Value = -2396.6667
$ws.Cells.Item(22, 14).Value = -2823.2222

$ws.Cells.Item(27, 8).Value = 2416.6
$ws.Cells.Item(27, 9).Value = 2691.6667
$ws.Cells.Item(27, 10).Value = 2233.2222
$ws.Cells.Item(27, 11).Value = 2691.6667
$ws.Cells.Item(27, 12).Value = 2233.2222
$ws.Cells.Item(27, 13).Value = -2584.6667
$ws.Cells.Item(27, 14).Value = -2447.2222

$ws.Cells.Item(40, 8).Value = 3977.6155
$ws.Cells.Item(40, 9).Value = 2458
$ws.Cells.Item(40, 11).Value = 2458
$ws.Cells.Item(40, 13).Value = -2322

$ws.Cells.Item(132, 8).Value = 2177.7083
$ws.Cells.Item(132, 9).Value = 1345.75
$ws.Cells.Item(132, 10).Value = 3009.6667
$ws.Cells.Item(132, 11).Value = 4037.25
$ws.Cells.Item(132, 12).Value = 9029.000100000001
$ws.Cells.Item(132, 13).Value = -1507.25
$ws.Cells.Item(132, 14).Value = -14089.0001

$ws.Cells.Item(136, 8).Value = 2810.6428
$ws.Cells.Item(136, 9).Value = 1427.2106
$ws.Cells.Item(136, 11).Value = 4281.6318
$ws.Cells.Item(136, 13).Value = -1731.6318

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(123, 8).Value = 47599.332
$ws.Cells.Item(123, 10).Value = 47599.332
$ws.Cells.Item(123, 12).Value = 47599.332
$ws.Cells.Item(123, 14).Value = -57399.332

$ws.Cells.Item(126, 8).Value = 6516
$ws.Cells.Item(126, 9).Value = 7251.647
$ws.Cells.Item(126, 10).Value = 5265.4
$ws.Cells.Item(126, 11).Value = 21754.941
$ws.Cells.Item(126, 12).Value = 15796.2
$ws.Cells.Item(126, 13).Value = -19284.941
$ws.Cells.Item(126, 14).Value = -20736.2

$ws.Cells.Item(132, 8).Value = 2560.25
$ws.Cells.Item(132, 9).Value = 1747.3334
$ws.Cells.Item(132, 11).Value = 5242.0002
$ws.Cells.Item(132, 13).Value = -2712.0002

$ws.Cells.Item(136, 8).Value = 12627802
$ws.Cells.Item(136, 9).Value = 15874271
$ws.Cells.Item(136, 10).Value = 2644.4443
$ws.Cells.Item(136, 11).Value = 47622813
$ws.Cells.Item(136, 12).Value = 7933.3329
$ws.Cells.Item(136, 13).Value = -47620263
$ws.Cells.Item(136, 14).Value = -13033.3329

Write-Output "Applied scheduled-runner updates to Sheets workbook."
